# Add the missing "DEFAULT" column (H) to the TABELS schema sheet, with a
# default value of 0 for the two INT11 fields that had one in the source
# edit (Partnership / PostCat), and move the selection to H20 to match the
# author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "DEFAULT" header cell next to TYPE/CO columns.
$ws.Range("H1").Value = "DEFAULT"

# Default values for the Partnership (row 9) and PostCat (row 22) fields.
$ws.Range("H9").Value = 0
$ws.Range("H22").Value = 0

# Best-effort: scroll the view so column B is the left-most visible column
# (matches topLeftCell="B1" in the target sheetView) before moving the
# selection to its final resting place.
try {
    $excel.ActiveWindow.ScrollColumn = 2
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}

# Final selection/active cell, matching the author's last edit location.
$ws.Range("H20").Select()
